# 23 dec 2023 update
# Record a new weekly payment entry (row 11) on the "MD10000.20-OCT" sheet:
#   B11 = 22-Dec-2023 (payment date)
#   C11 = 700         (amount paid)
#   D11 = 1           (payment count)
# The dependent totals (K1, O1, K2, O2) recalculate automatically from
# their existing formulas once the new values are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MD10000.20-OCT")
$ws.Activate()

$ws.Range("B11").Value = 45282   # 22-Dec-2023 (date serial number)
$ws.Range("C11").Value = 700
$ws.Range("D11").Value = 1

# Leave the active selection on D11, matching the saved cursor position.
$ws.Range("D11").Select()
